$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / tab
$ws.Name = "3_prospec.csv"

# Remove the old "in:vendor" / "ProSpec" column (old column C) entirely -
# Excel shifts D:K left into C:J, carrying values + styles with it.
$ws.Columns("C:C").Delete()

# Insert a brand-new first column (becomes the new "in:vendor_set" column),
# pushing the former A:J right into B:K.
$ws.Columns("A:A").Insert()
$ws.Columns("A:A").ColumnWidth = 13

# Header + data for the new column A
$ws.Range("A2").Value = "in:vendor_set"
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A3").Value = '$(ProSpec)'
$ws.Range("A4").Value = '$(ProSpec)'
$ws.Range("A5").Value = '$(ProSpec)'
$ws.Range("A6").Value = '$(ProSpec)'
$ws.Range("A7").Value = '$(ProSpec)'

# Update the selection to match the new layout
$ws.Range("A3:A7").Select()

Write-Output "done"
